$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172 (pushes existing rows 172-191 down to 173-192)
$ws.Rows("172:172").Insert()

# Populate the newly inserted row 172 with the new weekly price record
$ws.Cells.Item(172, 1).Value = 10
$ws.Cells.Item(172, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(172, 3).Value = "La Araucanía"
$ws.Cells.Item(172, 4).Value = 44449
$ws.Cells.Item(172, 5).Value = 9
$ws.Cells.Item(172, 6).Value = 100112037
$ws.Cells.Item(172, 7).Value = "Cebollín"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 115
$ws.Cells.Item(172, 11).Value = 8000
$ws.Cells.Item(172, 12).Value = 8000
$ws.Cells.Item(172, 13).Value = 8000
$ws.Cells.Item(172, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(172, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(172, 16).Value = 667
$ws.Cells.Item(172, 17).Value = 12
$ws.Cells.Item(172, 18).Value = "Hortaliza"
